# Settings back to qfin22 presentation - recalibration
$wb = $excel.ActiveWorkbook

# --- Update the long autocorrelation arrays (stored as text strings) ---
$wsAR = $wb.Worksheets.Item("AR")
$wsAR.Range("B5").Value = "[1.0, 0.22119845019596843, 0.16497097833142768, 0.15080475712961955, 0.12202387130058637, 0.3723644190998144, 0.15890356125730234, 0.16373383831362237, 0.1275904329900922, 0.10564284908993798, 0.025734576114490437, 0.06292062386278043, 0.06517541820195633, 0.09038811306262425, 0.12866920415251715, 0.04581300154645229, 0.05155633657867453, 0.05225625328876189, 0.04578381346466241, 0.08153853734744927]"

$wsSETAR = $wb.Worksheets.Item("SETAR")
$wsSETAR.Range("B10").Value = "[1.0, 0.25670036475624286, 0.20061771024202968, 0.1763131710456594, 0.1432118423813015, 0.28928909397178865, 0.14904541542512667, 0.1201673806150889, 0.10130692676865051, 0.11650153694056052, 0.045283374824621955, 0.06601077582992537, 0.042530789710169506, 0.06400302716491664, 0.09164318948053046, 0.03729027077774172, 0.050158944829187084, 0.03862347048177345, 0.033254331856330165, 0.04756200414655771]"

$wsGARCH = $wb.Worksheets.Item("GARCH")
$wsGARCH.Range("B6").Value = "[1.0, 0.003776397430326913, -0.0038506268908314475, -0.029931271341614954, -0.05216805559011897, 0.33590278377431443, 0.006456781155535773, 0.04057323911177375, -0.0034665417867512757, -0.020930902647498237, -0.05068775802995781, -0.048797830398633044, -0.0110381460612157, 0.018591745787028897, 0.05082791755580215, 0.002070425527879282, -0.011667002886938272, -0.030749915619392536, -0.006201038742565208, 0.02530631160773021]"

$wsTARCH = $wb.Worksheets.Item("TARCH")
$wsTARCH.Range("B6").Value = "[1.0, 0.0029429573660402493, -0.0065372242278268836, -0.03315628797347469, -0.052006541208077364, 0.3352527663776792, 0.00861113387387093, 0.043140721912395485, -0.0031371387485181726, -0.0199646907429981, -0.05122818922043748, -0.04780391242730359, -0.009959195787159132, 0.017574760158549325, 0.050720350208469, 0.0007950648799055855, -0.01224595525602734, -0.03182553065802484, -0.007064241286314997, 0.025457744342421238]"

$wsAR_TARCH = $wb.Worksheets.Item("AR_TARCH")
$wsAR_TARCH.Range("B6").Value = "[1.0, 0.03149255790457128, -0.0026330204502176074, -0.0007057522728670013, -0.03534504348194452, 0.24679879315407569, 0.002848195996758502, 0.03394824653373111, -0.00612784893800172, 0.0054189544663430095, -0.04965787879358057, -0.03930729007108875, -0.031182947527808442, 0.010258572554141131, 0.0555640381615112, -0.015339827641469649, -0.011805099828895552, -0.002983177433960842, -0.0026048921480626304, 0.021868926156090446]"

# --- Update numeric parameter values ---

# AR sheet
$wsAR.Range("B2").Value = 0.6291419451161907
$wsAR.Range("B3").Value = 0.8336228321218863
$wsAR.Range("B4").Value = 2679.221736904591

# SETAR sheet
$wsSETAR.Range("B4").Value = -16.95772502174763
$wsSETAR.Range("B5").Value = 0.7419219118194756
$wsSETAR.Range("B6").Value = 1880.479076701475
$wsSETAR.Range("B7").Value = 16.09120455330741
$wsSETAR.Range("B8").Value = 0.7575022577691661
$wsSETAR.Range("B9").Value = 1995.327320139124

# GARCH sheet
$wsGARCH.Range("B2").Value = 0.4271512959024379
$wsGARCH.Range("B3").Value = 129.1297637706128
$wsGARCH.Range("B4").Value = 0.1015001136922193
$wsGARCH.Range("B5").Value = 0.8544702078291411

# TARCH sheet
$wsTARCH.Range("B2").Value = 1.097871747516779
$wsTARCH.Range("B3").Value = 127.378381157966
$wsTARCH.Range("B4").Value = 0.1255879616819686
$wsTARCH.Range("B5").Value = 0.8552317805600218
$wsTARCH.Range("B7").Value = -0.04735204262921223

# AR_TARCH sheet
$wsAR_TARCH.Range("B2").Value = 1.525129102881922
$wsAR_TARCH.Range("B3").Value = 150.9070011067791
$wsAR_TARCH.Range("B4").Value = 0.139601241445151
$wsAR_TARCH.Range("B5").Value = 0.8284034503296247
$wsAR_TARCH.Range("B7").Value = -0.05038983331393765
$wsAR_TARCH.Range("B9").Value = 0.8131076963132683
